# Merge the split "<id>...</id>" runs (opening tag / value / closing tag,
# three separate <w:r> elements) into a single run for each of the three
# plain "p062v_N" identifiers (fig_p062v_1 / fig_p062v_2 are left as-is).
$d = $word.ActiveDocument

$ids = @("p062v_1", "p062v_2", "p062v_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $old, 2)
}
